$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Signed Ints" table (rows 14-15, mirroring the DEC table at rows 6-7) ---

# Copy the bordered/filled 2x3 table format from B6:D7 down to B14:D15
$ws.Range("B6:D7").Copy($ws.Range("B14"))

# Header row (row 14)
$ws.Range("B14").Value = "DECIMAL "
$ws.Range("C14").Value = "U_BINARY (up to +255)"
$ws.Range("D14").Value = "S_BINARY (-128 to +127)"

# Data row (row 15)
$ws.Range("B15").Value = 86
$ws.Range("C15").Formula = "=DEC2BIN(B15,8)"
$ws.Range("D15").Formula = "=BASE(MOD((B15*-1), 2^8), 2, 8)"

# Match row heights used throughout the sheet (14.4pt)
$ws.Rows(13).RowHeight = 14.4
$ws.Rows(14).RowHeight = 14.4
$ws.Rows(15).RowHeight = 14.4

# New E15 cell just to the right of the table, centered, no border/fill
$ws.Range("E15").HorizontalAlignment = -4108

# Column width adjustments to fit the new, longer header text
$ws.Columns("C").ColumnWidth = 20.49
$ws.Columns("D").ColumnWidth = 23.49
$ws.Columns("E").ColumnWidth = 29.79

# Update selection like the source workbook
$ws.Range("D15").Select()
